# Update currency year to 2023
$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsRACP  = $wb.Worksheets.Item("RACP")

# Update the shared string text "2021 dollars per 2012 dollar" -> "2023 dollars per 2012 dollar"
# (referenced from About!B11)
$wsAbout.Range("B11").Value = "2023 dollars per 2012 dollar"

# Update the computed conversion factor in About!A11
$wsAbout.Range("A11").Value = 0.75350342301658668

# Move the active selection on the About sheet to B12
$wsAbout.Activate()
$wsAbout.Range("B12").Select()

# Update the Cost Cap value on the RACP sheet from 160 to 180
$wsRACP.Range("B2").Value = 180
